$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188, shifting existing rows 188:309 down to 189:310
$ws.Rows(188).Insert()

# Populate the newly inserted row 188 with the new data point
$ws.Range("A188").Value = 10
$ws.Range("B188").Value = "Vega Modelo de Temuco"
$ws.Range("C188").Value = "La Araucanía"
$ws.Range("D188").Value = 44719
$ws.Range("E188").Value = 9
$ws.Range("F188").Value = 100112017
$ws.Range("G188").Value = "Apio"
$ws.Range("H188").Value = "Americana (o)"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 70
$ws.Range("K188").Value = 9000
$ws.Range("L188").Value = 10000
$ws.Range("M188").Value = 9429
$ws.Range("N188").Value = "$/docena de matas"
$ws.Range("O188").Value = "Provincia del Elquí"
$ws.Range("P188").Value = 1572
$ws.Range("Q188").Value = 6
$ws.Range("R188").Value = "Hortaliza"
